# Apply scheduled runner updates to market/profit data cells across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 464.84848
$ws.Range("I41").Value = 568.2727
$ws.Range("J41").Value = 413.13635
$ws.Range("K41").Value = 568.2727
$ws.Range("L41").Value = 413.13635
$ws.Range("M41").Value = -128.2727
$ws.Range("N41").Value = -1293.13635
$ws.Range("H70").Value = 972.94446
$ws.Range("I70").Value = 1154.2858
$ws.Range("J70").Value = 857.5454999999999
$ws.Range("K70").Value = 3462.8574
$ws.Range("L70").Value = 2572.6365
$ws.Range("M70").Value = -3192.8574
$ws.Range("N70").Value = -3112.6365
$ws.Range("H73").Value = 972.94446
$ws.Range("I73").Value = 1154.2858
$ws.Range("J73").Value = 857.5454999999999
$ws.Range("K73").Value = 3462.8574
$ws.Range("L73").Value = 2572.6365
$ws.Range("M73").Value = -2526.8574
$ws.Range("N73").Value = -4444.6365
$ws.Range("H132").Value = 5106679.5
$ws.Range("I132").Value = 5323965
$ws.Range("J132").Value = 475
$ws.Range("K132").Value = 15971895
$ws.Range("L132").Value = 1425
$ws.Range("M132").Value = -15969365
$ws.Range("N132").Value = -6485
$ws.Range("H135").Value = 1185.2122
$ws.Range("I135").Value = 699.1429000000001
$ws.Range("J135").Value = 3907.2
$ws.Range("K135").Value = 6292.2861
$ws.Range("L135").Value = 35164.8
$ws.Range("M135").Value = -3757.2861
$ws.Range("N135").Value = -40234.8
$ws.Range("H137").Value = 1839.95
$ws.Range("I137").Value = 1337.5
$ws.Range("J137").Value = 2174.9167
$ws.Range("K137").Value = 4012.5
$ws.Range("L137").Value = 6524.750100000001
$ws.Range("M137").Value = -1462.5
$ws.Range("N137").Value = -11624.7501
$ws.Range("H141").Value = 1989.0303
$ws.Range("I141").Value = 1907.75
$ws.Range("J141").Value = 4590
$ws.Range("K141").Value = 5723.25
$ws.Range("L141").Value = 13770
$ws.Range("M141").Value = -543.25
$ws.Range("N141").Value = -24130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 17520.4
$ws.Range("I6").Value = 50002
$ws.Range("J6").Value = 9400
$ws.Range("K6").Value = 50002
$ws.Range("L6").Value = 9400
$ws.Range("M6").Value = -49829
$ws.Range("N6").Value = -9746
$ws.Range("H32").Value = 19235.725
$ws.Range("I32").Value = 3539.4175
$ws.Range("J32").Value = 223287.72
$ws.Range("K32").Value = 3539.4175
$ws.Range("L32").Value = 223287.72
$ws.Range("M32").Value = -3252.4175
$ws.Range("N32").Value = -223861.72
$ws.Range("H61").Value = 1482.2543
$ws.Range("I61").Value = 909.7646999999999
$ws.Range("K61").Value = 909.7646999999999
$ws.Range("M61").Value = -697.7646999999999
$ws.Range("H63").Value = 2303.25
$ws.Range("I63").Value = 1297.3334
$ws.Range("K63").Value = 1297.3334
$ws.Range("M63").Value = -611.3334
$ws.Range("H66").Value = 2303.25
$ws.Range("I66").Value = 1297.3334
$ws.Range("K66").Value = 6486.666999999999
$ws.Range("M66").Value = -3054.666999999999
$ws.Range("H74").Value = 702.9375
$ws.Range("I74").Value = 646
$ws.Range("K74").Value = 646
$ws.Range("M74").Value = 228
$ws.Range("H77").Value = 702.9375
$ws.Range("I77").Value = 646
$ws.Range("K77").Value = 3230
$ws.Range("M77").Value = 1138
$ws.Range("H132").Value = 2066.6543
$ws.Range("I132").Value = 2098.446
$ws.Range("J132").Value = 1730.5714
$ws.Range("K132").Value = 6295.338
$ws.Range("L132").Value = 5191.7142
$ws.Range("M132").Value = -3765.338
$ws.Range("N132").Value = -10251.7142
$ws.Range("H136").Value = 1482.2543
$ws.Range("I136").Value = 909.7646999999999
$ws.Range("K136").Value = 2729.2941
$ws.Range("M136").Value = -179.2941000000001
$ws.Range("H138").Value = 59478.6
$ws.Range("J138").Value = 59478.6
$ws.Range("L138").Value = 59478.6
$ws.Range("N138").Value = -69758.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 29800
$ws.Range("J45").Value = 29800
$ws.Range("L45").Value = 29800
$ws.Range("N45").Value = -31416
$ws.Range("H56").Value = 25500
$ws.Range("J56").Value = 25500
$ws.Range("L56").Value = 25500
$ws.Range("N56").Value = -26978

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2750
$ws.Range("I62").Value = 2760
$ws.Range("K62").Value = 2760
$ws.Range("M62").Value = -2136
$ws.Range("H65").Value = 2750
$ws.Range("I65").Value = 2760
$ws.Range("K65").Value = 13800
$ws.Range("M65").Value = -10680
$ws.Range("H107").Value = 7642.8667
$ws.Range("I107").Value = 17839.834
$ws.Range("K107").Value = 17839.834
$ws.Range("M107").Value = -15919.834
$ws.Range("H132").Value = 38464428
$ws.Range("I132").Value = 40002950
$ws.Range("J132").Value = 35717056
$ws.Range("K132").Value = 120008850
$ws.Range("L132").Value = 107151168
$ws.Range("M132").Value = -120006320
$ws.Range("N132").Value = -107156228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1035.4225
$ws.Range("J131").Value = 998.32355
$ws.Range("L131").Value = 2994.97065
$ws.Range("N131").Value = -13074.97065
$ws.Range("H132").Value = 2810.3635
$ws.Range("I132").Value = 1702
$ws.Range("J132").Value = 3056.6667
$ws.Range("K132").Value = 15318
$ws.Range("L132").Value = 27510.0003
$ws.Range("M132").Value = -12788
$ws.Range("N132").Value = -32570.0003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5185
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 6496
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 6496
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -6720
$ws.Range("H40").Value = 47624.273
$ws.Range("I40").Value = 251869.75
$ws.Range("J40").Value = 2236.389
$ws.Range("K40").Value = 251869.75
$ws.Range("L40").Value = 2236.389
$ws.Range("M40").Value = -251733.75
$ws.Range("N40").Value = -2508.389
$ws.Range("H68").Value = 4067.2727
$ws.Range("I68").Value = 1488
$ws.Range("J68").Value = 6216.6665
$ws.Range("K68").Value = 1488
$ws.Range("L68").Value = 6216.6665
$ws.Range("M68").Value = -739
$ws.Range("N68").Value = -7714.6665
$ws.Range("H71").Value = 4067.2727
$ws.Range("I71").Value = 1488
$ws.Range("J71").Value = 6216.6665
$ws.Range("K71").Value = 7440
$ws.Range("L71").Value = 31083.3325
$ws.Range("M71").Value = -3696
$ws.Range("N71").Value = -38571.3325
$ws.Range("H122").Value = 2101.5293
$ws.Range("I122").Value = 2126.25
$ws.Range("J122").Value = 2042.2
$ws.Range("K122").Value = 6378.75
$ws.Range("L122").Value = 6126.6
$ws.Range("M122").Value = -3928.75
$ws.Range("N122").Value = -11026.6
$ws.Range("H126").Value = 5185
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 6496
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 19488
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -24428
$ws.Range("H136").Value = 1015.75
$ws.Range("I136").Value = 972.55884
$ws.Range("J136").Value = 1750
$ws.Range("K136").Value = 2917.67652
$ws.Range("L136").Value = 5250
$ws.Range("M136").Value = -367.67652
$ws.Range("N136").Value = -10350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6947122
$ws.Range("J62").Value = 3014.2856
$ws.Range("L62").Value = 3014.2856
$ws.Range("N62").Value = -4262.2856
$ws.Range("H65").Value = 6947122
$ws.Range("J65").Value = 3014.2856
$ws.Range("L65").Value = 15071.428
$ws.Range("N65").Value = -21311.428
$ws.Range("H122").Value = 2050.5
$ws.Range("I122").Value = 1641.2
$ws.Range("J122").Value = 2608.6365
$ws.Range("K122").Value = 4923.6
$ws.Range("L122").Value = 7825.9095
$ws.Range("M122").Value = -2473.6
$ws.Range("N122").Value = -12725.9095
